$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.190810322761536
$ws.Range("B1").Value = 2.284543037414551
$ws.Range("C1").Value = 4.596956253051758
$ws.Range("D1").Value = 3.431880235671997
$ws.Range("E1").Value = 1.196469187736511
